$d = $word.ActiveDocument

$replacements = @(
    @("963÷7=137, 4", "708÷6=118, 0"),
    @("496÷4=124, 0", "555÷6=92, 3"),
    @("877÷9=97, 4",  "732÷5=146, 2"),
    @("708÷3=236, 0", "257÷3=85, 2"),
    @("316÷4=79, 0",  "463÷2=231, 1"),
    @("314÷8=39, 2",  "558÷7=79, 5"),
    @("855÷4=213, 3", "233÷3=77, 2"),
    @("703÷3=234, 1", "705÷7=100, 5"),
    @("128÷5=25, 3",  "748÷2=374, 0"),
    @("634÷4=158, 2", "694÷5=138, 4"),
    @("971÷8=121, 3", "423÷5=84, 3"),
    @("644÷8=80, 4",  "983÷2=491, 1"),
    @("380÷3=126, 2", "687÷8=85, 7"),
    @("934÷7=133, 3", "962÷4=240, 2"),
    @("476÷6=79, 2",  "999÷4=249, 3"),
    @("542÷8=67, 6",  "205÷2=102, 1"),
    @("312÷9=34, 6",  "611÷2=305, 1"),
    @("323÷9=35, 8",  "254÷3=84, 2"),
    @("494÷4=123, 2", "548÷3=182, 2"),
    @("272÷6=45, 2",  "758÷9=84, 2"),
    @("838÷5=167, 3", "510÷6=85, 0"),
    @("208÷3=69, 1",  "772÷5=154, 2"),
    @("662÷7=94, 4",  "730÷2=365, 0"),
    @("981÷2=490, 1", "724÷9=80, 4"),
    @("826÷6=137, 4", "324÷2=162, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $new, 2)
}
